# Apply "writing all variable costs in opex" changes to the
# "Coupling Parameters" sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Country changed from DE to NL
$ws.Range("B1").Value = "NL"

# End Year: 2050 -> 2025
$ws.Range("B4").Value = 2025

# maximum_investment_capacity_per_year: 10000000 -> 10000
$ws.Range("B13").Value = 10000

# fix_prices_to_2030 description text updated to reference 2030 data
$ws.Range("C15").Value = "for verification runs. Fix fuel prices, CO2 and electricity demand to 2030 data"

# yearly_CO2_prices: FALSE -> TRUE
$ws.Range("B17").Value = $true

# targetinvestment_per_year: 2000 -> 1000
$ws.Range("B21").Value = 1000

# New row: writeALLcostsinOPEX flag set to TRUE
$ws.Range("A23").Value = "writeALLcostsinOPEX"
$ws.Range("B23").Value = $true

# Update the active selection to reflect the author's final cursor position
$ws.Activate()
$ws.Range("C12").Select()
